$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 75.71429000000001
$ws.Range("I11").Value = 75.71429000000001
$ws.Range("K11").Value = 75.71429000000001
$ws.Range("M11").Value = 64.28570999999999
$ws.Range("H39").Value = 337.25
$ws.Range("I39").Value = 292.9
$ws.Range("K39").Value = 878.6999999999999
$ws.Range("M39").Value = -582.6999999999999
$ws.Range("H40").Value = 5769.25
$ws.Range("I40").Value = 3044.4443
$ws.Range("K40").Value = 3044.4443
$ws.Range("M40").Value = -2869.4443
$ws.Range("H70").Value = 93354.82000000001
$ws.Range("I70").Value = 800
$ws.Range("J70").Value = 102610.3
$ws.Range("K70").Value = 2400
$ws.Range("L70").Value = 307830.9
$ws.Range("M70").Value = -2130
$ws.Range("N70").Value = -308370.9
$ws.Range("H73").Value = 93354.82000000001
$ws.Range("I73").Value = 800
$ws.Range("J73").Value = 102610.3
$ws.Range("K73").Value = 2400
$ws.Range("L73").Value = 307830.9
$ws.Range("M73").Value = -1464
$ws.Range("N73").Value = -309702.9

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 39871.89
$ws.Range("I2").Value = 61234.883
$ws.Range("K2").Value = 61234.883
$ws.Range("M2").Value = -61121.883
$ws.Range("H5").Value = 421.18182
$ws.Range("I5").Value = 92.666664
$ws.Range("K5").Value = 92.666664
$ws.Range("M5").Value = 19.333336
$ws.Range("H61").Value = 3430.3914
$ws.Range("I61").Value = 2059.0588
$ws.Range("J61").Value = 7315.8335
$ws.Range("K61").Value = 2059.0588
$ws.Range("L61").Value = 7315.8335
$ws.Range("M61").Value = -1847.0588
$ws.Range("N61").Value = -7739.8335
$ws.Range("H74").Value = 12536.667
$ws.Range("I74").Value = 13591
$ws.Range("K74").Value = 13591
$ws.Range("M74").Value = -12717
$ws.Range("H77").Value = 12536.667
$ws.Range("I77").Value = 13591
$ws.Range("K77").Value = 67955
$ws.Range("M77").Value = -63587
$ws.Range("H97").Value = 2549.7144
$ws.Range("I97").Value = 2549.7144
$ws.Range("K97").Value = 2549.7144
$ws.Range("M97").Value = -2053.7144
$ws.Range("H102").Value = 5909.1665
$ws.Range("I102").Value = 1819.3334
$ws.Range("K102").Value = 1819.3334
$ws.Range("M102").Value = -197.3334
$ws.Range("H116").Value = 39871.89
$ws.Range("I116").Value = 61234.883
$ws.Range("K116").Value = 61234.883
$ws.Range("M116").Value = -58940.883
$ws.Range("H123").Value = 0
$ws.Range("J123").Value = 0
$ws.Range("L123").Value = ""
$ws.Range("N123").Value = ""
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = ""
$ws.Range("N131").Value = ""
$ws.Range("H132").Value = 2491.0278
$ws.Range("I132").Value = 828.29034
$ws.Range("J132").Value = 12800
$ws.Range("K132").Value = 2484.87102
$ws.Range("L132").Value = 38400
$ws.Range("M132").Value = 45.12897999999996
$ws.Range("N132").Value = -43460
$ws.Range("H136").Value = 3430.3914
$ws.Range("I136").Value = 2059.0588
$ws.Range("J136").Value = 7315.8335
$ws.Range("K136").Value = 6177.176399999999
$ws.Range("L136").Value = 21947.5005
$ws.Range("M136").Value = -3627.176399999999
$ws.Range("N136").Value = -27047.5005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 39871.89
$ws.Range("I3").Value = 61234.883
$ws.Range("K3").Value = 61234.883
$ws.Range("M3").Value = -61120.883
$ws.Range("H4").Value = 421.18182
$ws.Range("I4").Value = 92.666664
$ws.Range("K4").Value = 92.666664
$ws.Range("M4").Value = 22.333336
$ws.Range("H94").Value = 1533.2106
$ws.Range("I94").Value = 1142.3334
$ws.Range("K94").Value = 1142.3334
$ws.Range("M94").Value = -691.3334
$ws.Range("H134").Value = 2920.9756
$ws.Range("I134").Value = 1834.4706
$ws.Range("K134").Value = 5503.4118
$ws.Range("M134").Value = -2968.4118

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H15").Value = 1888.5555
$ws.Range("J15").Value = 2418.1667
$ws.Range("L15").Value = 2418.1667
$ws.Range("N15").Value = -2758.1667
$ws.Range("H58").Value = 225155.67
$ws.Range("I58").Value = 386780.53
$ws.Range("J58").Value = 3984.7896
$ws.Range("K58").Value = 386780.53
$ws.Range("L58").Value = 3984.7896
$ws.Range("M58").Value = -386577.53
$ws.Range("N58").Value = -4390.7896
$ws.Range("H132").Value = 2412.7026
$ws.Range("I132").Value = 1653.9678
$ws.Range("J132").Value = 6332.8335
$ws.Range("K132").Value = 4961.903399999999
$ws.Range("L132").Value = 18998.5005
$ws.Range("M132").Value = -2431.903399999999
$ws.Range("N132").Value = -24058.5005
$ws.Range("H134").Value = 4450.706
$ws.Range("J134").Value = 7322.6665
$ws.Range("L134").Value = 21967.9995
$ws.Range("N134").Value = -27037.9995
$ws.Range("H136").Value = 225155.67
$ws.Range("I136").Value = 386780.53
$ws.Range("J136").Value = 3984.7896
$ws.Range("K136").Value = 1160341.59
$ws.Range("L136").Value = 11954.3688
$ws.Range("M136").Value = -1157791.59
$ws.Range("N136").Value = -17054.3688

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H99").Value = 3675
$ws.Range("I99").Value = 1025
$ws.Range("J99").Value = 5000
$ws.Range("K99").Value = 3075
$ws.Range("L99").Value = 15000
$ws.Range("M99").Value = -829
$ws.Range("N99").Value = -19492
$ws.Range("H132").Value = 2092.96
$ws.Range("I132").Value = 771
$ws.Range("J132").Value = 4443.1113
$ws.Range("K132").Value = 6939
$ws.Range("L132").Value = 39988.00169999999
$ws.Range("M132").Value = -4409
$ws.Range("N132").Value = -45048.00169999999
$ws.Range("H137").Value = 2881.625
$ws.Range("I137").Value = 2881.625
$ws.Range("K137").Value = 8644.875
$ws.Range("M137").Value = -3544.875

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H13").Value = 1911.7778
$ws.Range("I13").Value = 2300
$ws.Range("K13").Value = 2300
$ws.Range("M13").Value = -2161
$ws.Range("H62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = ""
$ws.Range("N62").Value = ""
$ws.Range("H65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = ""
$ws.Range("N65").Value = ""
$ws.Range("H97").Value = 893.0526
$ws.Range("I97").Value = 1030.3077
$ws.Range("J97").Value = 595.6667
$ws.Range("K97").Value = 1030.3077
$ws.Range("L97").Value = 595.6667
$ws.Range("M97").Value = -534.3077000000001
$ws.Range("N97").Value = -1587.6667
$ws.Range("H113").Value = 1003818.6
$ws.Range("I113").Value = 1253781.4
$ws.Range("K113").Value = 1253781.4
$ws.Range("M113").Value = -1251611.4
$ws.Range("H122").Value = 76427.28999999999
$ws.Range("J122").Value = 2400
$ws.Range("L122").Value = 7200
$ws.Range("N122").Value = -12100
$ws.Range("H128").Value = 78249
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 78249
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = ""
$ws.Range("M128").Value = ""
$ws.Range("N128").Value = -88209
$ws.Range("H129").Value = 75737
$ws.Range("J129").Value = 75737
$ws.Range("L129").Value = 75737
$ws.Range("N129").Value = -85737
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = ""
$ws.Range("N130").Value = ""
$ws.Range("H132").Value = 210817.77
$ws.Range("I132").Value = 258585.95
$ws.Range("K132").Value = 775757.8500000001
$ws.Range("M132").Value = -773227.8500000001

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 406228.16
$ws.Range("I7").Value = 630769.2
$ws.Range("K7").Value = 630769.2
$ws.Range("M7").Value = -630657.2
$ws.Range("H40").Value = 593907.8
$ws.Range("I40").Value = 774472.4399999999
$ws.Range("J40").Value = 7072.75
$ws.Range("K40").Value = 774472.4399999999
$ws.Range("L40").Value = 7072.75
$ws.Range("M40").Value = -774336.4399999999
$ws.Range("N40").Value = -7344.75
$ws.Range("H121").Value = 68248.92999999999
$ws.Range("J121").Value = 68248.92999999999
$ws.Range("L121").Value = 68248.92999999999
$ws.Range("N121").Value = -71742.92999999999
$ws.Range("H126").Value = 406228.16
$ws.Range("I126").Value = 630769.2
$ws.Range("K126").Value = 1892307.6
$ws.Range("M126").Value = -1889837.6
$ws.Range("H132").Value = 2898.6667
$ws.Range("I132").Value = 1983.8334
$ws.Range("J132").Value = 3630.5334
$ws.Range("K132").Value = 5951.5002
$ws.Range("L132").Value = 10891.6002
$ws.Range("M132").Value = -3421.5002
$ws.Range("N132").Value = -15951.6002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 15490
$ws.Range("J32").Value = 12000
$ws.Range("L32").Value = 12000
$ws.Range("N32").Value = -12634
$ws.Range("H95").Value = 82344
$ws.Range("J95").Value = 82344
$ws.Range("L95").Value = 82344
$ws.Range("N95").Value = -87836
$ws.Range("H113").Value = 1957.1333
$ws.Range("I113").Value = 2041.3334
$ws.Range("J113").Value = 1830.8334
$ws.Range("K113").Value = 6124.0002
$ws.Range("L113").Value = 5492.5002
$ws.Range("M113").Value = -3954.0002
$ws.Range("N113").Value = -9832.5002
$ws.Range("H132").Value = 1681.7
$ws.Range("I132").Value = 1822.5883
$ws.Range("J132").Value = 883.3333
$ws.Range("K132").Value = 5467.7649
$ws.Range("L132").Value = 2649.9999
$ws.Range("M132").Value = -2937.7649
$ws.Range("N132").Value = -7709.9999
$ws.Range("H136").Value = 347952.12
$ws.Range("I136").Value = 456885.1
$ws.Range("J136").Value = 5591.4287
$ws.Range("K136").Value = 1370655.3
$ws.Range("L136").Value = 16774.2861
$ws.Range("M136").Value = -1368105.3
$ws.Range("N136").Value = -21874.2861
